$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: id=3, nome=Magdo, valor=3, data=2025-06-01 ---
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Magdo"
$ws.Range("C4").Value = 3

# Copy D2 ("2025-06-01") into D4 so the value/type is reused verbatim
# (avoids Excel re-parsing the text as a date), then re-apply the
# standard data-row format from D3 on top of it.
$ws.Range("D2").Copy()
$ws.Range("D4").PasteSpecial(-4104)
$ws.Range("D3").Copy()
$ws.Range("D4").PasteSpecial(-4122)

# Match the formatting (borders/alignment) of the rest of the data rows
$ws.Range("A3:C3").Copy()
$ws.Range("A4:C4").PasteSpecial(-4122)

# --- Row 5: id=4, nome=Victor, valor=2, data=2025-06-01 ---
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Victor"
$ws.Range("C5").Value = 2

$ws.Range("D2").Copy()
$ws.Range("D5").PasteSpecial(-4104)
$ws.Range("D3").Copy()
$ws.Range("D5").PasteSpecial(-4122)

$ws.Range("A3:C3").Copy()
$ws.Range("A5:C5").PasteSpecial(-4122)

# --- H4: replicate the empty, underline-styled marker cell from H14 ---
$ws.Range("H14").Copy()
$ws.Range("H4").PasteSpecial(-4122)

# --- Update the active selection to H4 ---
$ws.Range("H4").Select()
